$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3885.5293
$ws.Range("I98").Value = 1843.5
$ws.Range("J98").Value = 13415
$ws.Range("K98").Value = 1843.5
$ws.Range("L98").Value = 13415
$ws.Range("M98").Value = -345.5
$ws.Range("N98").Value = -16411
$ws.Range("H122").Value = 3885.5293
$ws.Range("I122").Value = 1843.5
$ws.Range("J122").Value = 13415
$ws.Range("K122").Value = 5530.5
$ws.Range("L122").Value = 40245
$ws.Range("M122").Value = -3080.5
$ws.Range("N122").Value = -45145
$ws.Range("H135").Value = 4049.28
$ws.Range("I135").Value = 3892.3
$ws.Range("J135").Value = 4677.2
$ws.Range("K135").Value = 35030.7
$ws.Range("L135").Value = 42094.8
$ws.Range("M135").Value = -32495.7
$ws.Range("N135").Value = -47164.8
$ws.Range("H138").Value = 35062.742
$ws.Range("I138").Value = 2214.9443
$ws.Range("J138").Value = 80544.30499999999
$ws.Range("K138").Value = 6644.8329
$ws.Range("L138").Value = 241632.915
$ws.Range("M138").Value = -1504.8329
$ws.Range("N138").Value = -251912.915

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 21717.375
$ws.Range("I43").Value = 23657.25
$ws.Range("J43").Value = 19777.5
$ws.Range("K43").Value = 23657.25
$ws.Range("L43").Value = 19777.5
$ws.Range("M43").Value = -23344.25
$ws.Range("N43").Value = -20403.5
$ws.Range("H74").Value = 679915.1
$ws.Range("I74").Value = 1201281.6
$ws.Range("J74").Value = 28207
$ws.Range("K74").Value = 1201281.6
$ws.Range("L74").Value = 28207
$ws.Range("M74").Value = -1200407.6
$ws.Range("H77").Value = 679915.1
$ws.Range("I77").Value = 1201281.6
$ws.Range("J77").Value = 28207
$ws.Range("K77").Value = 6006408
$ws.Range("L77").Value = 141035
$ws.Range("M77").Value = -6002040
$ws.Range("H104").Value = 28999.4
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 28999.4
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 28999.4
$ws.Range("N104").Value = -35987.4
$ws.Range("H132").Value = 1481.0588
$ws.Range("I132").Value = 1296.5172
$ws.Range("J132").Value = 2551.4
$ws.Range("K132").Value = 3889.5516
$ws.Range("L132").Value = 7654.200000000001
$ws.Range("M132").Value = -1359.5516

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1469.85
$ws.Range("I94").Value = 655.5
$ws.Range("J94").Value = 4727.25
$ws.Range("K94").Value = 655.5
$ws.Range("L94").Value = 4727.25
$ws.Range("M94").Value = -204.5
$ws.Range("N94").Value = -5629.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 950
$ws.Range("I2").Value = 950
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 950
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -837
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1226
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H17").Value = 30000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 30000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30348
$ws.Range("H94").Value = 3449.7693
$ws.Range("I94").Value = 2071.2
$ws.Range("J94").Value = 4311.375
$ws.Range("K94").Value = 2071.2
$ws.Range("L94").Value = 4311.375
$ws.Range("M94").Value = -1620.2
$ws.Range("N94").Value = -5213.375
$ws.Range("H99").Value = 3489.4614
$ws.Range("I99").Value = 2033.091
$ws.Range("J99").Value = 11499.5
$ws.Range("K99").Value = 2033.091
$ws.Range("L99").Value = 11499.5
$ws.Range("M99").Value = -535.0909999999999
$ws.Range("H107").Value = 1462.5
$ws.Range("I107").Value = 1758.6666
$ws.Range("J107").Value = 1284.8
$ws.Range("K107").Value = 1758.6666
$ws.Range("L107").Value = 1284.8
$ws.Range("M107").Value = 161.3334
$ws.Range("H126").Value = 3489.4614
$ws.Range("I126").Value = 2033.091
$ws.Range("J126").Value = 11499.5
$ws.Range("K126").Value = 6099.272999999999
$ws.Range("L126").Value = 34498.5
$ws.Range("M126").Value = -3629.272999999999
$ws.Range("H132").Value = 34298
$ws.Range("I132").Value = 35308.277
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 105924.831
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -103394.831

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1164.8235
$ws.Range("I92").Value = 1300.2307
$ws.Range("J92").Value = 724.75
$ws.Range("K92").Value = 3900.6921
$ws.Range("L92").Value = 2174.25
$ws.Range("M92").Value = -2652.6921
$ws.Range("N92").Value = -4670.25
$ws.Range("H113").Value = 986.8889
$ws.Range("I113").Value = 599.6667
$ws.Range("J113").Value = 1180.5
$ws.Range("K113").Value = 1799.0001
$ws.Range("L113").Value = 3541.5
$ws.Range("M113").Value = 370.9999
$ws.Range("N113").Value = -7881.5
$ws.Range("H122").Value = 1207.0667
$ws.Range("I122").Value = 615.4
$ws.Range("J122").Value = 1502.9
$ws.Range("K122").Value = 5538.599999999999
$ws.Range("L122").Value = 13526.1
$ws.Range("M122").Value = -3088.599999999999
$ws.Range("H124").Value = 2950
$ws.Range("I124").Value = 2950
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 8850
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -3940
$ws.Range("H129").Value = 3280.7778
$ws.Range("I129").Value = 2121.2222
$ws.Range("J129").Value = 4440.3335
$ws.Range("K129").Value = 6363.6666
$ws.Range("L129").Value = 13321.0005
$ws.Range("M129").Value = -1363.6666
$ws.Range("N129").Value = -23321.0005
$ws.Range("H136").Value = 1586.5
$ws.Range("I136").Value = 1840
$ws.Range("J136").Value = 1333
$ws.Range("K136").Value = 5520
$ws.Range("L136").Value = 3999
$ws.Range("M136").Value = -420
$ws.Range("N136").Value = -14199

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 52000
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 100000
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 100000
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -100312
$ws.Range("H70").Value = 6313.5186
$ws.Range("I70").Value = 7071.4707
$ws.Range("J70").Value = 5025
$ws.Range("K70").Value = 7071.4707
$ws.Range("L70").Value = 5025
$ws.Range("M70").Value = -6801.4707
$ws.Range("H73").Value = 6313.5186
$ws.Range("I73").Value = 7071.4707
$ws.Range("J73").Value = 5025
$ws.Range("K73").Value = 7071.4707
$ws.Range("L73").Value = 5025
$ws.Range("M73").Value = -6135.4707
$ws.Range("H122").Value = 5029.933
$ws.Range("I122").Value = 4710
$ws.Range("J122").Value = 5669.8
$ws.Range("K122").Value = 14130
$ws.Range("L122").Value = 17009.4
$ws.Range("M122").Value = -11680
$ws.Range("N122").Value = -21909.4
$ws.Range("H133").Value = 106000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 106000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 106000
$ws.Range("M133").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 12669.667
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 12669.667
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 12669.667
$ws.Range("N14").Value = -13013.667
$ws.Range("H16").Value = 1636.2727
$ws.Range("I16").Value = 1149.9445
$ws.Range("J16").Value = 3824.75
$ws.Range("K16").Value = 1149.9445
$ws.Range("L16").Value = 3824.75
$ws.Range("M16").Value = -979.9445000000001
$ws.Range("H82").Value = 2170.4443
$ws.Range("I82").Value = 1933.7142
$ws.Range("J82").Value = 2999
$ws.Range("K82").Value = 1933.7142
$ws.Range("L82").Value = 2999
$ws.Range("M82").Value = -1572.7142
$ws.Range("H85").Value = 2170.4443
$ws.Range("I85").Value = 1933.7142
$ws.Range("J85").Value = 2999
$ws.Range("K85").Value = 1933.7142
$ws.Range("L85").Value = 2999
$ws.Range("M85").Value = -685.7141999999999
$ws.Range("H136").Value = 3745.7222
$ws.Range("I136").Value = 3192.9167
$ws.Range("J136").Value = 4851.3335
$ws.Range("K136").Value = 9578.750100000001
$ws.Range("L136").Value = 14554.0005
$ws.Range("M136").Value = -7028.750100000001
$ws.Range("H140").Value = 128800
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 128800
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 128800
$ws.Range("N140").Value = -139160

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 93249.836
$ws.Range("I3").Value = 173166.33
$ws.Range("J3").Value = 13333.333
$ws.Range("K3").Value = 173166.33
$ws.Range("L3").Value = 13333.333
$ws.Range("M3").Value = -173052.33
$ws.Range("H122").Value = 27149.52
$ws.Range("I122").Value = 31030.533
$ws.Range("J122").Value = 2200.1428
$ws.Range("K122").Value = 93091.599
$ws.Range("L122").Value = 6600.428400000001
$ws.Range("M122").Value = -90641.599
$ws.Range("H132").Value = 16803.588
$ws.Range("I132").Value = 20522.416
$ws.Range("J132").Value = 3415.8
$ws.Range("K132").Value = 61567.24800000001
$ws.Range("L132").Value = 10247.4
$ws.Range("M132").Value = -59037.24800000001
$ws.Range("H136").Value = 55271.582
$ws.Range("I136").Value = 105277.664
$ws.Range("J136").Value = 5265.5
$ws.Range("K136").Value = 315832.992
$ws.Range("L136").Value = 15796.5
$ws.Range("M136").Value = -313282.992
$ws.Range("N136").Value = -20896.5
